$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# "Table1" on the Data sheet currently covers EP 01..EP03 (A1:G4).
# The new episode, EP04, needs to be appended as a new table row with its
# Live / Views / Feedback / Discussions / Stars / Total Views metrics.
$tbl = $ws.ListObjects.Item("Table1")
$newRow = $tbl.ListRows.Add()

$newRow.Range.Item(1).Value = "EP04"
$newRow.Range.Item(2).Value = 4
$newRow.Range.Item(3).Value = 16
$newRow.Range.Item(4).Value = 0
$newRow.Range.Item(5).Value = 0
$newRow.Range.Item(6).Value = 2
$newRow.Range.Item(7).Value = 114

# Leave the cursor where the author last left it on the Data sheet.
$ws.Range("H7").Select() | Out-Null
